$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.075141191482544
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.153566479682922
